$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q4" sheet right before "2022-Q3", carrying
#    the per-fund holdings detail for the new quarter.
# ---------------------------------------------------------------------
$target = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($target)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'001305"
$q4.Range("C2").Value = "九泰天富改革新动力混合A"
$q4.Range("D2").Value = "'1.83"
$q4.Range("E2").Value = "'93.65"
$q4.Range("F2").Value = "'8.60"
$q4.Range("G2").Value = "'0.1574"
$q4.Range("H2").Value = 5

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'001844"
$q4.Range("C3").Value = "九泰久益灵活配置混合C"
$q4.Range("D3").Value = "'0.88"
$q4.Range("E3").Value = "'94.08"
$q4.Range("F3").Value = "'9.90"
$q4.Range("G3").Value = "'0.0871"
$q4.Range("H3").Value = 3

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'001782"
$q4.Range("C4").Value = "九泰久益灵活配置混合A"
$q4.Range("D4").Value = "'0.62"
$q4.Range("E4").Value = "'94.08"
$q4.Range("F4").Value = "'9.90"
$q4.Range("G4").Value = "'0.0614"
$q4.Range("H4").Value = 3

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'014938"
$q4.Range("C5").Value = "同泰产业升级混合A"
$q4.Range("D5").Value = "'0.95"
$q4.Range("E5").Value = "'68.98"
$q4.Range("F5").Value = "'2.89"
$q4.Range("G5").Value = "'0.0275"
$q4.Range("H5").Value = 8

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'009912"
$q4.Range("C6").Value = "九泰天富改革新动力混合C"
$q4.Range("D6").Value = "'0.03"
$q4.Range("E6").Value = "'93.65"
$q4.Range("F6").Value = "'8.60"
$q4.Range("G6").Value = "'0.0026"
$q4.Range("H6").Value = 5

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'014939"
$q4.Range("C7").Value = "同泰产业升级混合C"
$q4.Range("D7").Value = "'-0.01"
$q4.Range("E7").Value = "'68.98"
$q4.Range("F7").Value = "'2.89"
$q4.Range("G7").Value = "'-0.0003"
$q4.Range("H7").Value = 8

# Header + index cells carry the same bold/border style as the rest of the
# sheet (style already applied to A2 via the column-A border style); mirror
# the header row style from the sheet that sits right after it.
$q4.Range("B1:H1").Style = $target.Range("B1:H1").Style
$q4.Range("A2:A7").Style = $target.Range("A2:A2").Style

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4 and
#    push every following quarter down by one, re-numbering column A.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A2").EntireRow.Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.34
$summary.Range("A2").Style = $summary.Range("A3").Style

# Re-number column A (0-based position) for every row that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7
$summary.Range("A10").Value = 8
